# Financials update: insert a new "latest year" column before column D on the
# TLND sheet, pushing the existing D:K data right to E:L, then populate the
# brand-new column D with the new period's figures (period ending 2018-12-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D. Excel shifts D:K -> E:L and extends the sheet
#    dimension/used range to include the new column L automatically.
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D comes in with the left neighbour's (C's)
#    formatting. Copy the number formats/styles back from column E (which now
#    holds what used to be column D) so the new column matches the rest of
#    the data columns in each row (date style for row 7/38/80, number style
#    for the data rows, etc.)
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Match column D's width to the other data columns (E:I) so the new
#    column renders consistently with the rest of the block.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# 4) Fill in the new column D with the new period's values.

# 4a) "Period Ending" header rows (Income Statement / Balance Sheet / Cash
#     Flow Statement) get the new period date: 2018-12-31 (serial 43465).
$dateRows = @(7, 38, 80)
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 4).Value = 43465
}

# 4b) "Interest Expense" (row 22) has no figure for the new period -> "NA".
$ws.Cells.Item(22, 4).Value = "NA"

# 4c) All other data rows get their new-period numeric figure (including the
#     many rows whose figure is 0).
$numericMap = @{
    8=204300; 9=49500; 10=154800; 12=42400; 13=0; 14=0; 15=0; 17=245900;
    18=-41500; 20=900; 21=-36100; 23=-40700; 24=-300; 25=0; 26=-40400;
    27=-40400; 28=0; 29=0; 30=0; 31=0; 32=-900; 33=-40400; 34=0; 35=-40400;
    41=33700; 42=0; 43=78600; 44=0; 45=8300; 46=120700; 47=21600; 48=6300;
    49=69100; 50=0; 51=0; 52=1400; 53=0; 54=219100; 57=5800; 58=200;
    59=160900; 60=166900; 61=700; 62=27100; 63=0; 64=0; 65=0; 66=194700;
    68=0; 69=0; 70=0; 71=0; 72=-224200; 73=0; 74=0; 75=0; 76=24400; 77=0;
    81=-40400; 83=4600; 84=0; 85=0; 86=0; 87=0; 88=0; 89=3200; 91=-5000;
    92=0; 93=0; 94=-64500; 96=0; 97=0; 98=0; 99=0; 100=8600; 101=-600;
    102=-53300
}
foreach ($r in $numericMap.Keys) {
    $ws.Cells.Item($r, 4).Value = $numericMap[$r]
}
